# Fixed bugs in data
$wb = $excel.ActiveWorkbook

# --- Sessions sheet: correct the "Number of Redeems" (G) column values ---
$sessions = $wb.Worksheets.Item("Sessions")

$sessions.Range("G6").Value  = 2
$sessions.Range("G7").Value  = 2
$sessions.Range("G8").Value  = 2
$sessions.Range("G9").Value  = 2
$sessions.Range("G10").Value = 4
$sessions.Range("G12").Value = 4
$sessions.Range("G14").Value = 5
$sessions.Range("G15").Value = 16
$sessions.Range("G16").Value = 5
$sessions.Range("G17").Value = 5
$sessions.Range("G18").Value = 16
$sessions.Range("G19").Value = 16
$sessions.Range("G20").Value = 11
$sessions.Range("G21").Value = 11
$sessions.Range("G22").Value = 12
$sessions.Range("G23").Value = 12
$sessions.Range("G24").Value = 12
$sessions.Range("G25").Value = 13
$sessions.Range("G26").Value = 13
$sessions.Range("G27").Value = 13
$sessions.Range("G28").Value = 16
$sessions.Range("G29").Value = 2

# Highlight the corrected/flagged row in yellow
$sessions.Range("A34:H34").Interior.Color = 65535

# --- Restore cursor/selection positions left behind in various sheets ---
$employees = $wb.Worksheets.Item("Employees")
$employees.Range("F2").Select() | Out-Null

$specializes = $wb.Worksheets.Item("Specializes")
$specializes.Range("B3").Select() | Out-Null

$courses = $wb.Worksheets.Item("Courses")
$courses.Range("D11").Select() | Out-Null

$sessions.Activate() | Out-Null
$sessions.Range("M24").Select() | Out-Null
